$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 391; this pushes the existing rows 391..476 down to 392..477
# and extends the sheet dimension to A1:R477 automatically.
$ws.Rows("391:391").Insert()

# Populate the newly inserted row 391 with the new weekly record.
$ws.Cells.Item(391, 1).Value = 4
$ws.Cells.Item(391, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(391, 3).Value = "Los Lagos"
$ws.Cells.Item(391, 4).Value = 44964
$ws.Cells.Item(391, 5).Value = 10
$ws.Cells.Item(391, 6).Value = 100112023
$ws.Cells.Item(391, 7).Value = "Brócoli"
$ws.Cells.Item(391, 8).Value = "Sin especificar"
$ws.Cells.Item(391, 9).Value = "Primera"
$ws.Cells.Item(391, 10).Value = 1200
$ws.Cells.Item(391, 11).Value = 1500
$ws.Cells.Item(391, 12).Value = 1500
$ws.Cells.Item(391, 13).Value = 1500
$ws.Cells.Item(391, 14).Value = "$/unidad"
$ws.Cells.Item(391, 15).Value = "Región Metropolitana"
$ws.Cells.Item(391, 16).Value = 1500
$ws.Cells.Item(391, 17).Value = 1
$ws.Cells.Item(391, 18).Value = "Hortaliza"
